$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.033.24'

$ws.Range('E2').Value = '  +2.46%  '

$ws.Range('D3').Value = '3.212.80'

$ws.Range('E3').Value = '  +5.59%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.08'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  +3.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.72'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  +8.45%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  -0.36%  '

$ws.Range('D8').Value = '3.208.32'

$ws.Range('E8').Value = '  +5.81%  '

$ws.Range('E9').Value = '  +4.95%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.04'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  +9.84%  '

$ws.Range('E11').Value = '  +5.59%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.484'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  +5.58%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.16'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  +6.73%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000231'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  +5.85%  '

$ws.Range('D15').Value = '3.733.59'

$ws.Range('E15').Value = '  +5.93%  '

$ws.Range('D16').Value = '65.985.04'

$ws.Range('E16').Value = '  +2.31%  '

$ws.Range('D17').Value = '3.217.49'

$ws.Range('E17').Value = '  +5.51%  '

$ws.Range('B18').Value = 'TRON'

$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.114'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  +2.69%  '

$ws.Range('B19').Value = 'BitcoinCash'

$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '532.88'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  +10.72%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.09'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').Value = '  +7.53%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.48'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  +6.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.740'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  +8.37%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.70'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  +8.65%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.44'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  +8.22%  '

$ws.Range('E25').Value = '  +3.22%  '

$ws.Range('E26').Value = '  -0.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.29'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  +20.46%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.93'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  +7.72%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.25'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  +8.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '27.37'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  +6.29%  '

$ws.Range('E31').Value = '  -0.13%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.72'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  +5.43%  '

$ws.Range('E33').Value = '  +5.65%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '556.40'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  -1.02%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.32'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  +7.84%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.58'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  +4.42%  '

$ws.Range('B37').Value = 'OKB'

$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.63'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  +4.52%  '

$ws.Range('B38').Value = 'VeChain'

$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0450'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  +9.71%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0854'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').Value = '  +7.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.128'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  +7.43%  '

$ws.Range('D41').Value = '3.207.89'

$ws.Range('E41').Value = '  +10.65%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.88'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  +5.09%  '

$ws.Range('E43').Value = '  +4.58%  '

$ws.Range('E44').Value = '  +16.86%  '

$ws.Range('E45').Value = '  +13.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.29'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  +6.59%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  +0.03%  '

$ws.Range('D48').Value = '0.0₃0549'

$ws.Range('E48').Value = '  +4.33%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.63'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  +4.61%  '

$ws.Range('E50').Value = '  +4.26%  '

$ws.Range('E51').Value = '  +8.40%  '
